# Added results for Support Vector Machine
#
# Fills in the previously-empty "Support Vector Machine" rows (row 7 and
# row 14) on both result sheets, and gives the computed Balanced Accuracy
# column (F) the same 6-decimal number format already used by the other
# rows' F column. On the "With LOSOCV" sheet, cell E7 (Recall) additionally
# needs a new 5-decimal number format that isn't used anywhere else yet.

$wb = $excel.ActiveWorkbook

$sheetWithoutLOSOCV = $wb.Worksheets.Item("Dep or Non-Dep Without LOSOCV")
$sheetWithLOSOCV    = $wb.Worksheets.Item("Dep or Non-Dep With LOSOCV")

# ----- "Dep or Non-Dep Without LOSOCV" sheet -----

# Row 7: Support Vector Machine, with gender/age/work features
$sheetWithoutLOSOCV.Range("B7").Value = 0.85542168674698704
$sheetWithoutLOSOCV.Range("C7").Value = 0.86666666666666603
$sheetWithoutLOSOCV.Range("D7").Value = 0.76470588235294101
$sheetWithoutLOSOCV.Range("E7").Value = 0.8125
$sheetWithoutLOSOCV.Range("F7").Value = 0.84153661464585805
$sheetWithoutLOSOCV.Range("F7").NumberFormat = "0.000000"

# Row 14: Support Vector Machine, without gender/age/work features
$sheetWithoutLOSOCV.Range("B14").Value = 0.74698795180722799
$sheetWithoutLOSOCV.Range("C14").Value = 0.76
$sheetWithoutLOSOCV.Range("D14").Value = 0.55882352941176405
$sheetWithoutLOSOCV.Range("E14").Value = 0.644067796610169
$sheetWithoutLOSOCV.Range("F14").Value = 0.71818727490996404
$sheetWithoutLOSOCV.Range("F14").NumberFormat = "0.000000"

# ----- "Dep or Non-Dep With LOSOCV" sheet -----

# Row 7: Support Vector Machine, with gender/age/work features
$sheetWithLOSOCV.Range("B7").Value = 0.57816481498299599
$sheetWithLOSOCV.Range("C7").Value = 0.32727272727272699
$sheetWithLOSOCV.Range("D7").Value = 0.0775914994096812
$sheetWithLOSOCV.Range("E7").Value = 0.12232010413828499
$sheetWithLOSOCV.Range("E7").NumberFormat = "0.00000"
$sheetWithLOSOCV.Range("F7").Value = 0.57816481498299599
$sheetWithLOSOCV.Range("F7").NumberFormat = "0.000000"

# Row 14: Support Vector Machine, without gender/age/work features
$sheetWithLOSOCV.Range("B14").Value = 0.56467053148871305
$sheetWithLOSOCV.Range("C14").Value = 0.4
$sheetWithLOSOCV.Range("D14").Value = 0.122721369539551
$sheetWithLOSOCV.Range("E14").Value = 0.18257853257853199
$sheetWithLOSOCV.Range("F14").Value = 0.56467053148871305
$sheetWithLOSOCV.Range("F14").NumberFormat = "0.000000"
